# Updated docs and Portfolio Overview
# Updated documentation and Portfolio Overview narrative.
#
# Applies the narrative text edits on the architecture-diagram slide:
#   - "DATA SOURCES"            -> "MODERN DATA SOURCES"
#   - "Python ETL Scripts"      -> "Fabric Lakehouse Python Scripts for ETL"
#   - "SEMANTIC MODELING LAYER" -> "MODERN SEMANTIC MODELING LAYER"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$replacements = @(
    @{ Find = "DATA SOURCES"; Replace = "MODERN DATA SOURCES" },
    @{ Find = "Python ETL Scripts"; Replace = "Fabric Lakehouse Python Scripts for ETL" },
    @{ Find = "SEMANTIC MODELING LAYER"; Replace = "MODERN SEMANTIC MODELING LAYER" }
)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange

    foreach ($rep in $replacements) {
        # Skip work already done (also guards against double-application,
        # since some of the new strings contain the old string as a
        # substring, e.g. "MODERN SEMANTIC MODELING LAYER" contains
        # "SEMANTIC MODELING LAYER").
        if ($tr.Text -like ("*" + $rep.Replace + "*")) { continue }
        if ($tr.Text -notlike ("*" + $rep.Find + "*")) { continue }

        [void]$tr.Replace($rep.Find, $rep.Replace, 0, 0, 0)
    }
}
